$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in "Day 2" (column C) and "Day 3" (column D) sleep-diary answers for
# the last week block (rows 142-155), which were previously left blank.
# ---------------------------------------------------------------------------

# Row 142 - "您今天早上几点醒来?" (what time did you wake up?)
$ws.Range("C142").Value = 0.3125
$ws.Range("D142").Value = "8：51"

# Row 143 - "您今天几点起床?" (what time did you get out of bed?)
$ws.Range("C143").Value = 0.31944444444444442
$ws.Range("D143").Value = "9：00"

# Row 144 - "您昨晚几点上床?" (what time did you go to bed?)
$ws.Range("C144").Value = 0.95833333333333337
$ws.Range("D144").Value = "21：30"

# Row 145 - "您昨晚几点熄灯?" (what time did you turn off the lights?)
$ws.Range("C145").Value = 0.95833333333333337
$ws.Range("D145").Value = "21：30"

# Row 146 - minutes to fall asleep
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 0

# Row 147 - number of times woken up overnight
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 3

# Row 148 - total minutes awake overnight
$ws.Range("C148").Value = 5
$ws.Range("D148").Value = 30

# Row 149 - total minutes slept overnight
$ws.Range("C149").Value = 510
$ws.Range("D149").Value = 600

# Row 150 - used substances affecting sleep?
$ws.Range("C150").Value = "无"
$ws.Range("D150").Value = "无"

# Row 151 - used electronics before sleep / minutes
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 0

# Row 152-154 - rating scales
$ws.Range("C152").Value = 4
$ws.Range("D152").Value = 3

$ws.Range("C153").Value = 4
$ws.Range("D153").Value = 4

$ws.Range("C154").Value = 4
$ws.Range("D154").Value = 3

# Row 155 - nap during the day? / extra note
$ws.Range("C155").Value = "无"
$ws.Range("D155").Value = "无"
$ws.Range("E155").Value = "有 60 min"

# ---------------------------------------------------------------------------
# Move the selection to where the user left off editing.
# ---------------------------------------------------------------------------
$ws.Range("C157").Select()
